$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 updates
$ws.Range("A2").Value = "mTBhw539"
$ws.Range("B2").Value = 231006269
$ws.Range("C2").Value = "glzyafo15"
$ws.Range("D2").Value = "D7By#g9!"
$ws.Range("F2").Value = "gQhpCFEG"
$ws.Range("G2").Value = "HZCf"

# Row 3 updates
$ws.Range("A3").Value = "ltdqd443"
$ws.Range("B3").Value = 231006268
$ws.Range("C3").Value = "kjdkvhz38"
$ws.Range("D3").Value = "z7`$rU&6K"
$ws.Range("F3").Value = "xiMEYvUh"
$ws.Range("G3").Value = "plSY"
